# Dec 31 last updated
# Update the three tracker sheets (investment, growth, Expenses Ratio) with
# the latest entries for row 29 (month 12 / 2025) and move the active
# selection on each sheet as left by the author.

$wb = $excel.ActiveWorkbook

# --- investment sheet: fill in row 29 figures ---
$wsInvestment = $wb.Worksheets.Item("investment")
$wsInvestment.Range("C29").Value = 1103.5
$wsInvestment.Range("D29").Value = 6000
$wsInvestment.Range("E29").Value = 3000
$wsInvestment.Range("J29").Value = 310
$wsInvestment.Range("J30").Select()

# --- growth sheet: update row 29 figures ---
$wsGrowth = $wb.Worksheets.Item("growth")
$wsGrowth.Range("C29").Value = 6404.3
$wsGrowth.Range("D29").Value = 7340
$wsGrowth.Range("E29").Value = 6020
$wsGrowth.Range("J29").Value = 737.63
$wsGrowth.Range("C30").Select()

# --- Expenses Ratio sheet: just move the active selection ---
$wsExpenses = $wb.Worksheets.Item("Expenses Ratio")
$wsExpenses.Activate()
$wsExpenses.Range("G5").Select()
